# Applies the text edits described by the diff:
#  - Slide 6 "Content Placeholder 2": "Sub-TLVs Types:" ->
#        "Return Path Sub-TLVs Types:"
#  - Slide 7 "Content Placeholder 2": two sentence tweaks about the
#        Session-Reflector not transmitting / sending telemetry
#  - Slide 8 "Content Placeholder 2": "STAMP reply test packet" ->
#        "STAMP Session-Reflector reply test packet"
#  - Slide 9 "Content Placeholder 2": two sentence tweaks mentioning the
#        Session-Reflector / STAMP test session
#
# Helper functions read the shape's TextRange one character at a time
# (PowerPoint's TextRange.Characters(Start,Length) indexing reserves a
# slot - which reads back as an empty string - for each paragraph break,
# so a naive Python-style substring offset does not line up). Scanning
# the characters lets us locate the target phrase reliably regardless of
# how many paragraph breaks precede it, then grow the length to cover any
# such empty slots inside the match before rewriting just that sub-range.

function Get-FullCharList($tr) {
    $len = $tr.Length
    $chars = New-Object System.Collections.Generic.List[string]
    for ($i = 1; $i -le $len; $i++) {
        $chars.Add($tr.Characters($i, 1).Text)
    }
    return $chars
}

function Replace-InTextRange($tr, [string]$search, [string]$replace) {
    $chars = Get-FullCharList $tr
    $joined = [string]::Join("", $chars)
    $idx = $joined.IndexOf($search)
    if ($idx -lt 0) {
        throw "Substring not found: $search"
    }
    $length = $search.Length
    $blanksInside = 0
    for ($j = $idx; $j -lt ($idx + $length); $j++) {
        if ($chars[$j] -eq "") { $blanksInside++ }
    }
    $start = $idx + 1
    $adjLength = $length + $blanksInside
    $sub = $tr.Characters($start, $adjLength)
    if ($sub.Text -ne $search) {
        throw "Mismatch locating '$search': got '$($sub.Text)'"
    }
    $sub.Text = $replace
}

$p = $ppt.ActivePresentation

# Slide 6: "Sub-TLVs Types:" -> "Return Path Sub-TLVs Types:"
$sh = $p.Slides.Item(6).Shapes.Item(5)
Replace-InTextRange $sh.TextFrame.TextRange `
    "Sub-TLVs Types:" `
    "Return Path Sub-TLVs Types:"

# Slide 7: two sentence updates
$sh = $p.Slides.Item(7).Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Replace-InTextRange $tr `
    "The Session-Reflector does not transmit test packet back to the Session-Sender and terminates the test packet" `
    "The Session-Reflector does not transmit reply test packet to the Session-Sender and terminates the Session-Sender test packet"
Replace-InTextRange $tr `
    "Optionally, the Session-Reflector can send the performance metrics via streaming telemetry using the information from the received test packet" `
    "Optionally, the Session-Reflector can send the performance metrics via streaming telemetry using the information from the received Session-Sender test packet"

# Slide 8: one sentence update
$sh = $p.Slides.Item(8).Shapes.Item(2)
Replace-InTextRange $sh.TextFrame.TextRange `
    "The STAMP reply test packet may be transmitted to a different node than the Session-Sender " `
    "The STAMP Session-Reflector reply test packet may be transmitted to a different node than the Session-Sender "

# Slide 9: two sentence updates
$sh = $p.Slides.Item(9).Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Replace-InTextRange $tr `
    "For SR path, reply test packet may need to be sent in-band on a specific return SR path" `
    "For SR path, Session-Reflector reply test packet may need to be sent in-band on a specific return SR path"
Replace-InTextRange $tr `
    "Avoid signaling and maintaining dynamic state on Session-Reflector for the return path for each test session (each session-id, source-address) " `
    "Avoid signaling and maintaining dynamic state on Session-Reflector for the return path for each STAMP test session (each session-id, source-address) "
